# GPCRmd_B2AR_nomenclature_test.xlsx -- add more TM1 rows + shift TM separator rows up
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fill in the previously-blank rows 6-14 with the continuation of the TM1
#    helix table (generic numbering / sequence position / residue code).
# ---------------------------------------------------------------------------
$data = @(
    @(6,  "1x29", 1.39, "M40"),
    @(7,  "1x30", 1.40, "S41"),
    @(8,  "1x31", 1.41, "L42"),
    @(9,  "1x32", 1.42, "I43"),
    @(10, "1x33", 1.43, "V44"),
    @(11, "1x34", 1.44, "L45"),
    @(12, "1x35", 1.45, "A46"),
    @(13, "1x36", 1.46, "I47"),
    @(14, "1x37", 1.47, "V48")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Rows($r).RowHeight = 14
}

# ---------------------------------------------------------------------------
# 2) In each transmembrane block, the lone bold "separator" row (style of the
#    TM-header row) moves one row up - swap the formatting (and blank
#    contents) of the separator row with the plain row right before it.
# ---------------------------------------------------------------------------
$pairs = @(
    @(37, 38),
    @(69, 70),
    @(106, 107),
    @(134, 135),
    @(177, 178),
    @(216, 217),
    @(242, 243)
)

foreach ($p in $pairs) {
    $blankRow = $p[0]
    $sepRow = $p[1]
    $srcPlain = $blankRow - 1

    # Step 1: give A<blankRow> the bold "separator" formatting (copy it from
    # the separator row, which still has it at this point).
    $ws.Range("A$sepRow").Copy()
    $ws.Range("A$blankRow").PasteSpecial(-4122)

    # Step 2: give the whole old separator row the ordinary plain-row
    # formatting (copy it from a known plain row just above the block).
    $ws.Range("A$srcPlain" + ":C$srcPlain").Copy()
    $ws.Range("A$sepRow" + ":C$sepRow").PasteSpecial(-4122)

    # Step 3: drop B/C on the (now bold, single-cell) separator row - this
    # clears the clipboard, so it must be the last step for this block.
    $ws.Range("B$blankRow" + ":C$blankRow").Clear()
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Drop the now-unused last row (256) -- the sheet shrinks from 256 to 255
#    rows.
# ---------------------------------------------------------------------------
$ws.Rows(256).Delete()

# ---------------------------------------------------------------------------
# 4) Restore the active selection to A15 (matches the saved view state).
# ---------------------------------------------------------------------------
$ws.Range("A15").Select()
